$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50, shifting existing rows 50-118 down to 51-119
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with data
$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 44638
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100109
$ws.Range("H50").Value = "Uva"
$ws.Range("I50").Value = 100109001
$ws.Range("J50").Value = "Uva"
$ws.Range("K50").Value = "Red Globe"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 220
$ws.Range("N50").Value = 6500
$ws.Range("O50").Value = 7000
$ws.Range("P50").Value = 6773
$ws.Range("Q50").Value = "`$/caja 10 kilos"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 677
$ws.Range("T50").Value = 10
